# Swap the order of "dnasr281@gmail.com" and "System" in the
# "Recorded By" column (column G) wherever the old ordering is found.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Session Analysis Results")

$oldText = "dnasr281@gmail.com, System"
$newText = "System, dnasr281@gmail.com"

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count + $usedRange.Row - 1

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)  # Column G
    if ($cell.Value2 -eq $oldText) {
        $cell.Value2 = $newText
    }
}
